$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.170.17'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").Value = '1.904.85'
$ws.Range("E3").Value = '  +0.75%  '

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'306.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '

$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").Value = "'0.5237"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.78%  '

$ws.Range("D8").Value = "'0.3770"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.54%  '

$ws.Range("D9").Value = "'0.07249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").Value = "'21.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.04%  '

$ws.Range("D11").Value = "'0.9029"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("D12").Value = "'0.08530"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.58%  '

$ws.Range("D13").Value = '1.922.69'
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").Value = "'97.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.36%  '

$ws.Range("D15").Value = "'5.293"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").Value = "'1.0000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("D17").Value = "'0.000008636"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.90%  '

$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").Value = '27.208.77'
$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").Value = '2.153.54'
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("D24").Value = "'6.440"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.60%  '

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").Value = "'2.298"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.86%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'147.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'1.749"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.87%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = "'4.919"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.88%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'4.816"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.43%  '

$ws.Range("D32").Value = "'0.09288"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").Value = "'0.8056"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.94%  '

$ws.Range("D34").Value = "'0.05055"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.74%  '

$ws.Range("D35").Value = "'1.243"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.52%  '

$ws.Range("D36").Value = "'3.447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.84%  '

$ws.Range("D37").Value = "'2.957"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.40%  '

$ws.Range("D38").Value = "'2.619"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.41%  '

$ws.Range("D39").Value = "'0.5709"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.27%  '

$ws.Range("D40").Value = "'0.02000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("E41").Value = '  -0.12%  '

$ws.Range("D42").Value = "'9.160"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.89%  '

$ws.Range("D43").Value = "'6.638"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").Value = "'116.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.26%  '

$ws.Range("D45").Value = "'0.1518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").Value = "'0.4872"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.47%  '

$ws.Range("D47").Value = "'0.9996"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("D48").Value = "'10.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.35%  '

$ws.Range("D49").Value = "'1.617"
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("E51").Value = '  +0.25%  '
